$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows with both Price (D) and Volume(1h) (E) changes ---
$ws.Range("D2").Value = '58.284.91'
$ws.Range("E2").Value = '  -4.12%  '
$ws.Range("D3").Value = '2.644.18'
$ws.Range("E3").Value = '  -1.76%  '
$ws.Range("D5").Value = '''521.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.95%  '
$ws.Range("D6").Value = '''144.19'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.54%  '
$ws.Range("D8").Value = '''0.571'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.40%  '
$ws.Range("D9").Value = '''6.69'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.10%  '
$ws.Range("D11").Value = '''0.338'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.93%  '
$ws.Range("D12").Value = '''0.131'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.59%  '
$ws.Range("D13").Value = '3.108.50'
$ws.Range("E13").Value = '  -1.77%  '
$ws.Range("D14").Value = '58.290.51'
$ws.Range("E14").Value = '  -3.94%  '
$ws.Range("D15").Value = '''20.87'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.12%  '
$ws.Range("D17").Value = '2.654.87'
$ws.Range("E17").Value = '  -11.61%  '
$ws.Range("D20").Value = '''10.45'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.28%  '
$ws.Range("D21").Value = '''6.29'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.81%  '
$ws.Range("D22").Value = '''0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").Value = '''64.39'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.22%  '
$ws.Range("D24").Value = '''0.424'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.49%  '
$ws.Range("D26").Value = '''1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.65%  '
$ws.Range("D27").Value = '0.0₃0796'
$ws.Range("E27").Value = '  -2.81%  '
$ws.Range("D28").Value = '''7.11'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.15%  '
$ws.Range("D29").Value = '''6.63'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.16%  '
$ws.Range("D32").Value = '''152.97'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.06%  '
$ws.Range("D33").Value = '''18.82'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.64%  '
$ws.Range("D34").Value = '''4.13'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.08%  '
$ws.Range("D35").Value = '''1.18'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.95%  '
$ws.Range("D36").Value = '''0.905'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.10%  '
$ws.Range("D37").Value = '''0.859'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.59%  '
$ws.Range("D38").Value = '''36.82'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.42%  '
$ws.Range("D39").Value = '''1.44'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.09%  '
$ws.Range("D42").Value = '''0.607'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.74%  '
$ws.Range("D45").Value = '''19.38'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.10%  '
$ws.Range("D48").Value = '2.039.71'
$ws.Range("E48").Value = '  -4.58%  '
$ws.Range("D51").Value = '''18.27'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.95%  '

# --- Rows with only Volume(1h) (E) changes ---
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("E10").Value = '  -3.27%  '
$ws.Range("E16").Value = '  -1.53%  '
$ws.Range("E19").Value = '  -2.77%  '
$ws.Range("E25").Value = '  -2.11%  '
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("E40").Value = '  -0.93%  '
$ws.Range("E41").Value = '  +0.46%  '

# --- Rows with only Price (D) changes ---
$ws.Range("D18").Value = '''337.54'
$ws.Range("D18").Style = "Normal"

# --- Rows with full Coin/Link/Price/Volume swap (reordering) ---
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").Value = '''0.0970'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.40%  '
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").Value = '''268.91'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.60%  '
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").Value = '''0.0540'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("B47").Value = 'WhiteBITCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D47").Value = '''10.63'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.54%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '''4.67'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.28%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = '''0.0228'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.18%  '
